# Reset the deck's theme colour scheme from the custom "Integral" /
# "Red Violet" palette back to the stock Office colour scheme
# (Design > Variants > Colors > "Office" in the PowerPoint UI).
#
# The font scheme (Arial/Arial) and the format scheme (fills, lines,
# effects) are already identical between the two themes, so only the
# 12 theme colours need to change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Index -> (role, new RGB as 0xBBGGRR COM colour, i.e. R + G*256 + B*65536)
$colors.Item(1).RGB  = 0         # dk1      #000000
$colors.Item(2).RGB  = 16777215  # lt1      #FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      #44546A
$colors.Item(4).RGB  = 15132391  # lt2      #E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  #5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  #ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  #A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  #FFC000
$colors.Item(9).RGB  = 12874308  # accent5  #4472C4
$colors.Item(10).RGB = 4697456   # accent6  #70AD47
$colors.Item(11).RGB = 12673797  # hlink    #0563C1
$colors.Item(12).RGB = 7491477   # folHlink #954F72
